$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Frequência (batching moving)" table (A18:B22), mirroring the
# existing A10:F14 summary table's layout/style.
$ws.Range("A18").Value = "Projeto"
$ws.Range("B18").Value = "Frequência (batching moving)"

$ws.Range("A19").Value = "clojure"
$ws.Range("B19").Value = 0.25

$ws.Range("A20").Value = "gradle"
$ws.Range("B20").Value = 0.01

$ws.Range("A21").Value = "junit4"
$ws.Range("B21").Value = 0.06

$ws.Range("A22").Value = "spring-framework"
$ws.Range("B22").Value = 0.04

# Match the centered "GENERAL" style used by the other header/data blocks.
$ws.Range("A18:B22").HorizontalAlignment = -4108

# Column B now holds the long "Frequência (batching moving)" header, so it
# needs to widen; column A narrows slightly once it is split out on its own.
$ws.Range("A1").ColumnWidth = 14
$ws.Range("B1").ColumnWidth = 23.5
$ws.Range("C1").ColumnWidth = 17.83
$ws.Range("D1").ColumnWidth = 13.5
$ws.Range("E1").ColumnWidth = 12.5
$ws.Range("F1").ColumnWidth = 21.5

# Mirror where the author left the selection after entering the new table.
$ws.Range("A18:B22").Select()
